$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("O3").Value = 1.48

# Row 4
$ws.Range("G4").Value = 1.42
$ws.Range("H4").Value = 4.5
$ws.Range("I4").Value = 7
$ws.Range("J4").Value = 1.02
$ws.Range("K4").Value = 21
$ws.Range("L4").Value = 1.11
$ws.Range("M4").Value = 6.5
$ws.Range("N4").Value = 1.4
$ws.Range("O4").Value = 2.88
$ws.Range("P4").Value = 1.22
$ws.Range("Q4").Value = 4
$ws.Range("R4").Value = 1.57
$ws.Range("S4").Value = 2.25
$ws.Range("T4").Value = 11
$ws.Range("U4").Value = 9
$ws.Range("V4").Value = 8.5
$ws.Range("W4").Value = 11
$ws.Range("Z4").Value = 21
$ws.Range("AA4").Value = 9.5
$ws.Range("AB4").Value = 15
$ws.Range("AC4").Value = 41
$ws.Range("AI4").Value = 41
$ws.Range("AJ4").Value = 126

# Row 6
$ws.Range("H6").Value = 3.5
$ws.Range("K6").Value = 8.25
$ws.Range("L6").Value = 1.27
$ws.Range("M6").Value = 3.55
$ws.Range("N6").Value = 1.8
$ws.Range("O6").Value = 1.95
$ws.Range("P6").Value = 1.39
$ws.Range("Q6").Value = 2.85
$ws.Range("R6").Value = 1.65
$ws.Range("S6").Value = 2.1
$ws.Range("T6").Value = 8.25
$ws.Range("Y6").Value = 27
$ws.Range("Z6").Value = 8.25
$ws.Range("AA6").Value = 7.1
$ws.Range("AB6").Value = 14
$ws.Range("AC6").Value = 60
$ws.Range("AD6").Value = 10
$ws.Range("AE6").Value = 18
$ws.Range("AH6").Value = 28
$ws.Range("AI6").Value = 35
$ws.Range("AJ6").Value = 450

# Row 8
$ws.Range("G8").Value = 1.95
$ws.Range("H8").Value = 3.2
$ws.Range("I8").Value = 4.1
$ws.Range("J8").Value = 1.1
$ws.Range("K8").Value = 6.1
$ws.Range("L8").Value = 1.47
$ws.Range("M8").Value = 2.55
$ws.Range("N8").Value = 2.4
$ws.Range("O8").Value = 1.53
$ws.Range("P8").Value = 1.53
$ws.Range("Q8").Value = 2.4
$ws.Range("R8").Value = 2.1
$ws.Range("S8").Value = 1.65
$ws.Range("T8").Value = 5.6
$ws.Range("U8").Value = 8.5
$ws.Range("V8").Value = 9.5
$ws.Range("W8").Value = 18
$ws.Range("X8").Value = 20
$ws.Range("Y8").Value = 45
$ws.Range("Z8").Value = 6.1
$ws.Range("AA8").Value = 6.6
$ws.Range("AB8").Value = 21
$ws.Range("AD8").Value = 8.75
$ws.Range("AE8").Value = 22
$ws.Range("AG8").Value = 75
$ws.Range("AI8").Value = 70

# Row 11
$ws.Range("J11").Value = 1.1
$ws.Range("K11").Value = 7

# Row 12
$ws.Range("G12").Value = 1.85
$ws.Range("I12").Value = 4.1
$ws.Range("L12").Value = 1.37
$ws.Range("N12").Value = 2.07
$ws.Range("U12").Value = 7.9
$ws.Range("Z12").Value = 8.25
$ws.Range("AD12").Value = 10.25
$ws.Range("AE12").Value = 22
$ws.Range("AG12").Value = 65
$ws.Range("AI12").Value = 50

# Row 14
$ws.Range("J14").Value = 1.07
$ws.Range("K14").Value = 9
$ws.Range("N14").Value = 2.15
$ws.Range("O14").Value = 1.67

# Row 18
$ws.Range("N18").Value = 1.98
$ws.Range("O18").Value = 1.83

# Row 20
$ws.Range("G20").Value = 2.2
$ws.Range("I20").Value = 3.4
$ws.Range("J20").Value = 1.08
$ws.Range("K20").Value = 8
$ws.Range("N20").Value = 2.35
$ws.Range("O20").Value = 1.57
$ws.Range("AE20").Value = 15

# Row 22
$ws.Range("I22").Value = 4.25
$ws.Range("K22").Value = 6.9
$ws.Range("L22").Value = 1.34
$ws.Range("M22").Value = 3
$ws.Range("N22").Value = 2
$ws.Range("O22").Value = 1.72
$ws.Range("P22").Value = 1.45
$ws.Range("Q22").Value = 2.57
$ws.Range("R22").Value = 1.87
$ws.Range("S22").Value = 1.83
$ws.Range("T22").Value = 6.4
$ws.Range("U22").Value = 8
$ws.Range("W22").Value = 14.5
$ws.Range("Y22").Value = 29
$ws.Range("Z22").Value = 6.9
$ws.Range("AC22").Value = 80
$ws.Range("AF22").Value = 14
$ws.Range("AJ22").Value = 700

# Row 24
$ws.Range("L24").Value = 1.44
$ws.Range("M24").Value = 2.63

# Row 26
$ws.Range("G26").Value = 2.5
$ws.Range("I26").Value = 2.82
$ws.Range("N26").Value = 2.27
$ws.Range("P26").Value = 1.5
$ws.Range("Q26").Value = 2.27
$ws.Range("T26").Value = 6.7
$ws.Range("U26").Value = 11.25
$ws.Range("V26").Value = 10
$ws.Range("W26").Value = 27
$ws.Range("X26").Value = 24
$ws.Range("AA26").Value = 5.8
$ws.Range("AD26").Value = 7.1
$ws.Range("AE26").Value = 13
$ws.Range("AF26").Value = 10.75
$ws.Range("AG26").Value = 35
$ws.Range("AH26").Value = 29

# Row 28
$ws.Range("G28").Value = 2.7
$ws.Range("H28").Value = 3.1
$ws.Range("I28").Value = 2.52
$ws.Range("L28").Value = 1.45
$ws.Range("M28").Value = 2.37
$ws.Range("N28").Value = 2.32
$ws.Range("O28").Value = 1.47
$ws.Range("P28").Value = 1.5
$ws.Range("Q28").Value = 2.27
$ws.Range("R28").Value = 2
$ws.Range("U28").Value = 12
$ws.Range("V28").Value = 10.75
$ws.Range("W28").Value = 30
$ws.Range("X28").Value = 28
$ws.Range("Z28").Value = 7
$ws.Range("AB28").Value = 18.5
$ws.Range("AC28").Value = 120
$ws.Range("AD28").Value = 6.5
$ws.Range("AE28").Value = 11
$ws.Range("AF28").Value = 10.25
$ws.Range("AG28").Value = 27
$ws.Range("AH28").Value = 25

# Row 31
$ws.Range("G31").Value = 3.8
$ws.Range("J31").Value = 1.06
$ws.Range("K31").Value = 10
$ws.Range("L31").Value = 1.25
$ws.Range("M31").Value = 3.75
$ws.Range("N31").Value = 1.9
$ws.Range("O31").Value = 1.9
$ws.Range("P31").Value = 1.36
$ws.Range("Q31").Value = 3
$ws.Range("R31").Value = 1.75
$ws.Range("S31").Value = 2
$ws.Range("T31").Value = 12
$ws.Range("U31").Value = 19
$ws.Range("X31").Value = 29
$ws.Range("Y31").Value = 34
$ws.Range("Z31").Value = 11
$ws.Range("AB31").Value = 13
$ws.Range("AC31").Value = 41
$ws.Range("AD31").Value = 8
$ws.Range("AE31").Value = 10
$ws.Range("AH31").Value = 15
$ws.Range("AJ31").Value = 201

# Row 32
$ws.Range("L32").Value = 1.34
$ws.Range("M32").Value = 3
$ws.Range("P32").Value = 1.45
$ws.Range("Q32").Value = 2.55
$ws.Range("T32").Value = 6.7
$ws.Range("U32").Value = 9.25
$ws.Range("V32").Value = 8.75
$ws.Range("W32").Value = 18
$ws.Range("Y32").Value = 30
$ws.Range("AB32").Value = 14.5
$ws.Range("AC32").Value = 70
$ws.Range("AD32").Value = 10
$ws.Range("AE32").Value = 19
$ws.Range("AH32").Value = 32

Write-Host "Updated 190 cells across rows 3,4,6,8,11,12,14,18,20,22,24,26,28,31,32"
